# Änderung für Parallelisierung der Feldprüfung
# Appends 16 additional test rows (4 repeating blocks of the existing
# "BARVERKAUF 1" / "test2" / "BARVERKAUF" / "test" pattern) to Tabelle1
# and moves the active selection to D14, matching the author's edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(10027, 50000, "BARVERKAUF 1"),
    @(10027, 70003, "test2"),
    @(10026, 50000, "BARVERKAUF 1"),
    @(10026, 70003, "test2"),
    @(10028, 50000, "BARVERKAUF"),
    @(10028, 70003, "test"),
    @(10027, 50000, "BARVERKAUF 1"),
    @(10027, 70003, "test2"),
    @(10026, 50000, "BARVERKAUF 1"),
    @(10026, 70003, "test2"),
    @(10028, 50000, "BARVERKAUF"),
    @(10028, 70003, "test"),
    @(10027, 50000, "BARVERKAUF 1"),
    @(10027, 70003, "test2"),
    @(10026, 50000, "BARVERKAUF 1"),
    @(10026, 70003, "test2")
)

$row = 5
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row = $row + 1
}

$ws.Range("D14").Select()
